$d = $word.ActiveDocument

# 1) Merge "Data:3" + "1" + "/08/2025" runs into a single run "Data:31/08/2025"
$d.Content.Find.Execute("Data:31/08/2025", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Data:31/08/2025", 2) | Out-Null

# 2) Merge "Hora: " + "20:30" runs into a single run "Hora: 20:30"
$d.Content.Find.Execute("Hora: 20:30", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Hora: 20:30", 2) | Out-Null

# 3) Fill in the previously-empty participant e-mail cells in the table
$table = $d.Tables.Item(1)
$table.Cell(9, 2).Range.Text = "gustavofbraga@gmail.com"
$table.Cell(10, 2).Range.Text = "pedroffn209@gmail.com"
$table.Cell(12, 2).Range.Text = "vaasgalinari@sga.pucminas.br"
